$wb = $excel.ActiveWorkbook

# Mirror the user's navigation: select column D on the Debts sheet first
$wsDebts = $wb.Worksheets.Item("Debts")
$wsDebts.Activate() | Out-Null
$wsDebts.Columns("D").Select() | Out-Null

# Then go to the Fixed Assets sheet and insert a new column D (shifting
# the existing "basis"... columns to the right), mirroring the Debts layout
# which already has a "year" column in D.
$ws = $wb.Worksheets.Item("Fixed Assets")
$ws.Activate() | Out-Null
$ws.Columns("D").Select() | Out-Null
$ws.Columns("D").Insert() | Out-Null
$ws.Range("D1").Value = "year"
